$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("B4").Value = 0.5845681264288881
$ws.Range("C4").Value = 0.601
$ws.Range("D4").Value = 0.5729440557234217
$ws.Range("E4").Value = 0.5780000000000001
$ws.Range("F4").Value = 0.5717826952400543
$ws.Range("G4").Value = 0.572
$ws.Range("H4").Value = 0.576880256254212
$ws.Range("I4").Value = 0.575
$ws.Range("J4").Value = 0.5133491337347751
$ws.Range("K4").Value = 0.511
$ws.Range("L4").Value = 0.523151499444
$ws.Range("M4").Value = 0.522
$ws.Range("N4").Value = 0.6036753743931195
$ws.Range("O4").Value = 0.6159999999999999
$ws.Range("P4").Value = 0.6028518898370733
$ws.Range("Q4").Value = 0.599
$ws.Range("R4").Value = 0.5785652225833577
$ws.Range("S4").Value = 0.59
$ws.Range("T4").Value = 0.5729014797117272
$ws.Range("U4").Value = 0.5770000000000001
$ws.Range("V4").Value = 0.567879956246941
$ws.Range("W4").Value = 0.5670000000000001
$ws.Range("X4").Value = 0.5744155766747584
$ws.Range("Y4").Value = 0.5725
$ws.Range("Z4").Value = 0.6034130210693378
$ws.Range("AA4").Value = 0.615
$ws.Range("AB4").Value = 0.6028699011759008
$ws.Range("AC4").Value = 0.599

# Row 5
$ws.Range("B5").Value = 0.6548313376184899
$ws.Range("C5").Value = 0.835
$ws.Range("D5").Value = 0.5404544515196925
$ws.Range("E5").Value = 0.5620000000000001
$ws.Range("F5").Value = 0.4835430185074124
$ws.Range("G5").Value = 0.5279999999999999
$ws.Range("H5").Value = 0.5581453462618027
$ws.Range("I5").Value = 0.5325
$ws.Range("J5").Value = 0.677460317875371
$ws.Range("K5").Value = 0.9620000000000001
$ws.Range("L5").Value = 0.523779124346939
$ws.Range("M5").Value = 0.542
$ws.Range("N5").Value = 0.5238951285472165
$ws.Range("O5").Value = 0.585
$ws.Range("P5").Value = 0.5653831561873679
$ws.Range("Q5").Value = 0.5499999999999999
$ws.Range("R5").Value = 0.6613780097729983
$ws.Range("S5").Value = 0.8700000000000001
$ws.Range("T5").Value = 0.5350181940300203
$ws.Range("U5").Value = 0.5565
$ws.Range("V5").Value = 0.5085330649152446
$ws.Range("W5").Value = 0.5659999999999999
$ws.Range("X5").Value = 0.5662291051619348
$ws.Range("Y5").Value = 0.5445
$ws.Range("Z5").Value = 0.5418062130768486
$ws.Range("AA5").Value = 0.615
$ws.Range("AB5").Value = 0.5656486288646392
$ws.Range("AC5").Value = 0.554

# Row 6
$ws.Range("B6").Value = 0.5881671239646633
$ws.Range("C6").Value = 0.574
$ws.Range("D6").Value = 0.6217608926845342
$ws.Range("E6").Value = 0.6125
$ws.Range("F6").Value = 0.6805563181047538
$ws.Range("G6").Value = 0.669
$ws.Range("H6").Value = 0.7131871496881104
$ws.Range("I6").Value = 0.6930000000000001
$ws.Range("J6").Value = 0.543696350475867
$ws.Range("K6").Value = 0.541
$ws.Range("L6").Value = 0.5553221257548028
$ws.Range("M6").Value = 0.5525
$ws.Range("N6").Value = 0.6975852276227466
$ws.Range("O6").Value = 0.6849999999999999
$ws.Range("P6").Value = 0.7284479389363702
$ws.Range("Q6").Value = 0.7125000000000001
$ws.Range("R6").Value = 0.589092567594468
$ws.Range("S6").Value = 0.5740000000000001
$ws.Range("T6").Value = 0.6209315849744621
$ws.Range("U6").Value = 0.6129999999999999
$ws.Range("V6").Value = 0.6871353930080181
$ws.Range("W6").Value = 0.6679999999999999
$ws.Range("X6").Value = 0.7272190730520369
$ws.Range("Y6").Value = 0.7050000000000001
$ws.Range("Z6").Value = 0.7018982612123865
$ws.Range("AA6").Value = 0.6840000000000001
$ws.Range("AB6").Value = 0.7403448515741375
$ws.Range("AC6").Value = 0.721
